# Health and elemental mechanic now being used
# Update the second weakness-chart grid (rows 10-14) so that it keys off
# elemental names (Water, Earth, Lightning, Fire) instead of color names
# (Blue, Green, Yellow, Red), and move the active selection to H12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (row 10) - column labels for the matrix
# New shared strings must be introduced in "Earth" then "Water" order so
# they land at shared-string indices 8 and 9 respectively.
$ws.Range("G10").Value = "Earth"
$ws.Range("F10").Value = "Water"
$ws.Range("H10").Value = "Lightning"
$ws.Range("I10").Value = "Fire"

# Row labels (column E, rows 11-14)
$ws.Range("E12").Value = "Earth"
$ws.Range("E11").Value = "Water"
$ws.Range("E13").Value = "Lightning"
$ws.Range("E14").Value = "Fire"

# H11 was previously an implicit blank cell; now explicitly holds 0
$ws.Range("H11").Value = 0

# Move the selection/active cell to H12
$ws.Range("H12").Select()
